# Revert "Merge branch 'master' of https://github.com/vidwalk/SEP3"
#
# This script:
#  1. Strips the paragraph-mark run properties (w:pPr/w:rPr/w:lang) off the
#     "In order to make the system secure..." paragraph.
#  2. Inserts a new empty paragraph right after it (before the _GoBack
#     bookmark paragraph).
#  3. Deletes the paragraphs that discussed each individual threat
#     (Spoofing identity / Tampering / Information disclosure / Denial of
#     service / Elevation of privileges) which used to follow the
#     bookmark paragraph, right up to the section break.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: rebuild the "In order to make the system secure..." paragraph
# without its w:pPr (which only carried a w:rPr/w:lang="en-US").
# The paragraph is identified by its (unique) leading text.
# ---------------------------------------------------------------------
$introPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("In order to make the system secure")) {
        $introPara = $p
        break
    }
}

$introXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>In order to make the system secure</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">possible </w:t></w:r><w:r><w:t>risks and danger</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> that come with the project were discussed. All goals, means and powers of attacker were taken under consideration. Matters such as thread frequency and effect, preventive and corrective measures are examined in the next part of the document. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$introPara.Range.InsertXML($introXml)

# ---------------------------------------------------------------------
# Step 2: insert a brand-new empty paragraph directly before the
# paragraph that carries the _GoBack bookmark. The bookmark itself
# pinpoints exactly where that paragraph starts, regardless of
# paragraph index shifting from step 1.
# ---------------------------------------------------------------------
$goBackMark = $d.Bookmarks("_GoBack")
$goBackPara = $d.Range($goBackMark.Range.Start, $goBackMark.Range.Start).Paragraphs(1)
$goBackPara.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------
# Step 3: delete every paragraph from "The threat of spoofing identity..."
# through "...permit attackers to succeed in elevating their privileges."
# i.e. everything that used to follow the _GoBack bookmark paragraph,
# up to (but not including) the final section properties.
# ---------------------------------------------------------------------
$goBackMark = $d.Bookmarks("_GoBack")
$goBackPara = $d.Range($goBackMark.Range.Start, $goBackMark.Range.Start).Paragraphs(1)

$lastBodyPara = $d.Paragraphs($d.Paragraphs.Count)

if ($lastBodyPara.Range.Start -gt $goBackPara.Range.End) {
    $deleteRange = $d.Range($goBackPara.Range.End, $lastBodyPara.Range.End)
    $deleteRange.Delete()
}
